$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10:H11").NumberFormat = "@"

$ws.Range("A10").Value = "IEAGHG_sinter_plant"
$ws.Range("B10").Value = "sinter plant"
$ws.Range("C10").Value = "sinter"
$ws.Range("D10").Value = "output"
$ws.Range("E10").Value = "data/steel/SteelUnits_Variables.xlsx"
$ws.Range("F10").Value = "Sinter Plant"
$ws.Range("G10").Value = "data/steel/SteelUnits_Relationships.xlsx"
$ws.Range("H10").Value = "Sinter Plant"

$ws.Range("B11").Value = "blast furnace"
$ws.Range("C11").Value = "hot metal"
$ws.Range("D11").Value = "outflow"
$ws.Range("E11").Value = "data/steel/SteelUnits_Variables.xlsx"
$ws.Range("F11").Value = "Blast Furnace"
$ws.Range("G11").Value = "data/steel/SteelUnits_Relationships.xlsx"
$ws.Range("H11").Value = "Blast Furnace"
$ws.Range("A11").Value = "IEAGHG_blast_furnace"

$ws.Range("G15").Select()
